# Applies the "Rol" / "Operador" menu assignment change to the Hoja1 sheet.
$xlCenter = -4108  # xlCenter

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2: change F2 from "ADMIN" to "ADMIN_OPE" and add H2 = "Operador"
$ws.Range("F2").Value = "ADMIN_OPE"
$ws.Range("H2").Value = "Operador"
$ws.Range("H2").HorizontalAlignment = $xlCenter
$ws.Range("H2").VerticalAlignment = $xlCenter

# Row 3: change F3 from "ADMIN" to "ADMIN_OPE" and add H3 = "Operador"
$ws.Range("F3").Value = "ADMIN_OPE"
$ws.Range("H3").Value = "Operador"
$ws.Range("H3").HorizontalAlignment = $xlCenter
$ws.Range("H3").VerticalAlignment = $xlCenter

# Update the active selection to H2, matching the sheet view in the diff
$ws.Activate()
$ws.Range("H2").Select()
